$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(7)
$ws.Name = "shortest path length"

$ws.Range("B2").Value = 0.8758970319761407
$ws.Range("C2").Value = 0.9498616017991128
$ws.Range("D2").Value = 1.021659216609922
$ws.Range("E2").Value = 0.9575776925724472
$ws.Range("F2").Value = 1.123959771136158
$ws.Range("G2").Value = -45.29555860042941
$ws.Range("H2").Value = [double]"2.061914574464136E-244"
$ws.Range("B3").Value = 0.8723161531677792
$ws.Range("C3").Value = 0.951411570366714
$ws.Range("D3").Value = 1.026795463341972
$ws.Range("E3").Value = 0.953924457764149
$ws.Range("F3").Value = 0.9915445247219833
$ws.Range("G3").Value = -10.57115129478604
$ws.Range("H3").Value = [double]"7.846855882528151E-25"
$ws.Range("B4").Value = 0.8712117501784894
$ws.Range("C4").Value = 0.9458143118017956
$ws.Range("D4").Value = 1.028324662806458
$ws.Range("E4").Value = 0.9542580008471837
$ws.Range("F4").Value = 0.9083462792725192
$ws.Range("G4").Value = 13.03733874163368
$ws.Range("H4").Value = [double]"5.407789124429374E-36"
$ws.Range("B5").Value = 0.8063122115469168
$ws.Range("C5").Value = 0.9372692729918557
$ws.Range("D5").Value = 1.102544166829774
$ws.Range("E5").Value = 0.9674620621966543
$ws.Range("F5").Value = 0.9938308271423179
$ws.Range("G5").Value = -3.737391803586525
$ws.Range("H5").Value = 0.0001964838673362682
$ws.Range("B6").Value = 0.808670053995864
$ws.Range("C6").Value = 0.9377413113038149
$ws.Range("D6").Value = 1.095271585580679
$ws.Range("E6").Value = 0.9628120632249443
$ws.Range("F6").Value = 1.132941944872307
$ws.Range("G6").Value = -25.09273986495125
$ws.Range("H6").Value = [double]"3.830967716944697E-108"
$ws.Range("B7").Value = 0.73545178808999
$ws.Range("C7").Value = 0.9105364269970809
$ws.Range("D7").Value = 1.136409870733305
$ws.Range("E7").Value = 0.9632039921764856
$ws.Range("F7").Value = 0.973341474553722
$ws.Range("G7").Value = -1.030227299713375
$ws.Range("H7").Value = 0.3031526729863639
$ws.Range("B8").Value = 0.7762725124600383
$ws.Range("C8").Value = 0.9067769702746171
$ws.Range("D8").Value = 1.070299735906753
$ws.Range("E8").Value = 0.9427657846349355
$ws.Range("F8").Value = 1.031685828319533
$ws.Range("G8").Value = -11.98879126680685
$ws.Range("H8").Value = [double]"4.888628694190928E-31"
$ws.Range("A9").Value = 13
$ws.Range("B9").Value = 0.8063122115469168
$ws.Range("C9").Value = 0.9372692729918557
$ws.Range("D9").Value = 1.102544166829774
$ws.Range("E9").Value = 0.9674620621966543
$ws.Range("F9").Value = 1.147226879075929
$ws.Range("G9").Value = -25.47906792609648
$ws.Range("H9").Value = [double]"9.827984240912329E-111"
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 0.7536694317031035
$ws.Range("C10").Value = 0.9122959730744027
$ws.Range("D10").Value = 1.122635521534521
$ws.Range("E10").Value = 0.958319376005242
$ws.Range("F10").Value = 1.261240543277536
$ws.Range("G10").Value = -34.33595117226105
$ws.Range("H10").Value = [double]"2.898937046884094E-171"
$ws.Range("B11").Value = 0.8023392150341201
$ws.Range("C11").Value = 0.9399580445847646
$ws.Range("D11").Value = 1.087759542767883
$ws.Range("E11").Value = 0.9611093168159223
$ws.Range("F11").Value = 1.191832469569707
$ws.Range("G11").Value = -34.54622199356913
$ws.Range("H11").Value = [double]"1.054805864712953E-172"
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 0.7536694317031035
$ws.Range("C12").Value = 0.9122959730744027
$ws.Range("D12").Value = 1.122635521534521
$ws.Range("E12").Value = 0.958319376005242
$ws.Range("F12").Value = 1.243260558765712
$ws.Range("G12").Value = -32.297929610958
$ws.Range("H12").Value = [double]"2.760737292678883E-157"
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = 0.7536694317031035
$ws.Range("C13").Value = 0.9122959730744027
$ws.Range("D13").Value = 1.122635521534521
$ws.Range("E13").Value = 0.958319376005242
$ws.Range("F13").Value = 1.264753082742598
$ws.Range("G13").Value = -34.73409562894935
$ws.Range("H13").Value = [double]"5.470883975330541E-174"
$ws.Range("B14").Value = 0.8063122115469168
$ws.Range("C14").Value = 0.9372692729918557
$ws.Range("D14").Value = 1.102544166829774
$ws.Range("E14").Value = 0.9674620621966543
$ws.Range("F14").Value = 0.5948406788699871
$ws.Range("G14").Value = 52.81370237687918
$ws.Range("H14").Value = [double]"2.084995407862569E-291"
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 0.6210579615178411
$ws.Range("C15").Value = 0.8367531423125035
$ws.Range("D15").Value = 1.106650979924942
$ws.Range("E15").Value = 0.953769247036397
$ws.Range("F15").Value = 0.5389736105319591
$ws.Range("G15").Value = 31.10798395075452
$ws.Range("H15").Value = [double]"4.127015914543601E-149"
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = 0.6210579615178411
$ws.Range("C16").Value = 0.8367531423125035
$ws.Range("D16").Value = 1.106650979924942
$ws.Range("E16").Value = 0.953769247036397
$ws.Range("F16").Value = 0.5518139036158931
$ws.Range("G16").Value = 30.1450142470611
$ws.Range("H16").Value = [double]"1.691621094726068E-142"
$ws.Range("A17").Value = 24
$ws.Range("B17").Value = 0.808670053995864
$ws.Range("C17").Value = 0.9377413113038149
$ws.Range("D17").Value = 1.095271585580679
$ws.Range("E17").Value = 0.9628120632249443
$ws.Range("F17").Value = 0.7899824694278736
$ws.Range("G17").Value = 25.49092491055818
$ws.Range("H17").Value = [double]"8.181142713052413E-111"
$ws.Range("A18").Value = 23
$ws.Range("B18").Value = 0.808670053995864
$ws.Range("C18").Value = 0.9377413113038149
$ws.Range("D18").Value = 1.095271585580679
$ws.Range("E18").Value = 0.9628120632249443
$ws.Range("F18").Value = 0.8138135164115088
$ws.Range("G18").Value = 21.97604406258778
$ws.Range("H18").Value = [double]"1.249498057837654E-87"
$ws.Range("B19").Value = 0.808670053995864
$ws.Range("C19").Value = 0.9377413113038149
$ws.Range("D19").Value = 1.095271585580679
$ws.Range("E19").Value = 0.9628120632249443
$ws.Range("F19").Value = 0.9403321234969235
$ws.Range("G19").Value = 3.315603786430721
$ws.Range("H19").Value = 0.0009473617293182079
$ws.Range("I19").Value = "Yes. Right."
$ws.Range("B20").Value = 0.808670053995864
$ws.Range("C20").Value = 0.9377413113038149
$ws.Range("D20").Value = 1.095271585580679
$ws.Range("E20").Value = 0.9628120632249443
$ws.Range("F20").Value = 1.196209758401484
$ws.Range("G20").Value = -34.42421515512096
$ws.Range("H20").Value = [double]"7.212373963431037E-172"
$ws.Range("B21").Value = 0.73545178808999
$ws.Range("C21").Value = 0.9105364269970809
$ws.Range("D21").Value = 1.136409870733305
$ws.Range("E21").Value = 0.9632039921764856
$ws.Range("F21").Value = 1.814358332836917
$ws.Range("G21").Value = -86.49903451244867
$ws.Range("B22").Value = 0.7025284942713607
$ws.Range("C22").Value = 0.8957397589404499
$ws.Range("D22").Value = 1.146511685629753
$ws.Range("E22").Value = 0.961509989873116
$ws.Range("F22").Value = 1.263529318321717
$ws.Range("G22").Value = -26.48211863876961
$ws.Range("H22").Value = [double]"1.693780842806571E-117"
$ws.Range("B23").Value = 0.73545178808999
$ws.Range("C23").Value = 0.9105364269970809
$ws.Range("D23").Value = 1.136409870733305
$ws.Range("E23").Value = 0.9632039921764856
$ws.Range("F23").Value = 1.506960247016541
$ws.Range("G23").Value = -55.25953262163311
$ws.Range("H23").Value = [double]"4.81362194782138E-306"
$ws.Range("B24").Value = 0.8023392150341201
$ws.Range("C24").Value = 0.9399580445847646
$ws.Range("D24").Value = 1.087759542767883
$ws.Range("E24").Value = 0.9611093168159223
$ws.Range("F24").Value = 0.9477264980176456
$ws.Range("G24").Value = 2.003812030075474
$ws.Range("H24").Value = 0.04535983184228964
$ws.Range("I24").Value = "Yes. Right."
$ws.Range("B25").Value = 0.7762725124600383
$ws.Range("C25").Value = 0.9067769702746171
$ws.Range("D25").Value = 1.070299735906753
$ws.Range("E25").Value = 0.9427657846349355
$ws.Range("F25").Value = 0.8459855840829211
$ws.Range("G25").Value = 13.04854985556862
$ws.Range("H25").Value = [double]"4.76914511517546E-36"
$ws.Range("B26").Value = 0.7025284942713607
$ws.Range("C26").Value = 0.8957397589404499
$ws.Range("D26").Value = 1.146511685629753
$ws.Range("E26").Value = 0.961509989873116
$ws.Range("F26").Value = 0.9343579357042153
$ws.Range("G26").Value = 2.380787757792448
$ws.Range("H26").Value = 0.01746230149665006
$ws.Range("B27").Value = 0.808670053995864
$ws.Range("C27").Value = 0.9377413113038149
$ws.Range("D27").Value = 1.095271585580679
$ws.Range("E27").Value = 0.9628120632249443
$ws.Range("F27").Value = 0.9890605187687749
$ws.Range("G27").Value = -3.871428466536448
$ws.Range("H27").Value = 0.0001152289023049078
$ws.Range("B28").Value = 0.8063122115469168
$ws.Range("C28").Value = 0.9372692729918557
$ws.Range("D28").Value = 1.102544166829774
$ws.Range("E28").Value = 0.9674620621966543
$ws.Range("F28").Value = 0.9312971772374882
$ws.Range("G28").Value = 5.125850410611105
$ws.Range("H28").Value = [double]"3.555279154994583E-07"
$ws.Range("B29").Value = 0.8138290316641124
$ws.Range("C29").Value = 0.9257411037761902
$ws.Range("D29").Value = 1.065827737372534
$ws.Range("E29").Value = 0.9469957709699602
$ws.Range("F29").Value = 0.9290659438851491
$ws.Range("G29").Value = 3.166774833611174
$ws.Range("H29").Value = 0.00158811648486056
$ws.Range("I29").Value = "Yes. Right."
$ws.Range("B30").Value = 0.7536694317031035
$ws.Range("C30").Value = 0.9122959730744027
$ws.Range("D30").Value = 1.122635521534521
$ws.Range("E30").Value = 0.958319376005242
$ws.Range("F30").Value = 1.241419252993388
$ws.Range("G30").Value = -32.08921859329948
$ws.Range("H30").Value = [double]"7.492326899039872E-156"
$ws.Range("B31").Value = 0.8224276744631035
$ws.Range("C31").Value = 0.9310822452799874
$ws.Range("D31").Value = 1.056768156636358
$ws.Range("E31").Value = 0.9503786635943035
$ws.Range("F31").Value = 1.290935355821875
$ws.Range("G31").Value = -61.65479599148845
$ws.Range("B32").Value = 0.8063122115469168
$ws.Range("C32").Value = 0.9372692729918557
$ws.Range("D32").Value = 1.102544166829774
$ws.Range("E32").Value = 0.9674620621966543
$ws.Range("F32").Value = 0.8494522405376039
$ws.Range("G32").Value = 16.72618877372848
$ws.Range("H32").Value = [double]"1.483991031145629E-55"
$ws.Range("B33").Value = 0.8224276744631035
$ws.Range("C33").Value = 0.9310822452799874
$ws.Range("D33").Value = 1.056768156636358
$ws.Range("E33").Value = 0.9503786635943035
$ws.Range("F33").Value = 1.094809340200118
$ws.Range("G33").Value = -26.14790460524454
$ws.Range("H33").Value = [double]"3.075653947999153E-115"
$ws.Range("A34").Value = 14
$ws.Range("B34").Value = 0.8063122115469168
$ws.Range("C34").Value = 0.9372692729918557
$ws.Range("D34").Value = 1.102544166829774
$ws.Range("E34").Value = 0.9674620621966543
$ws.Range("F34").Value = 1.11646810031572
$ws.Range("G34").Value = -21.11945503320518
$ws.Range("H34").Value = [double]"3.807140826399548E-82"
